$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndexes = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$rowValues = @(1.02, 1.029535611330794, 1.0378792253427, 1.039116639903256, 1.049503707582211, 1.033788878336445, 1.034681847077829, 1.04066885129537, 1.041902741802839, 1.052260583338637, 1.015570603495936)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(2, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.030388896884912, 1.038528353158067, 1.03986801930943, 1.050324086768794, 1.033921319834343, 1.035176697750805, 1.041128302905858, 1.042464428939411, 1.052893200797284, 1.015734712811213)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(3, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.03094149738673, 1.038948381565166, 1.040354971912009, 1.050855589642355, 1.034005242070008, 1.035496720919522, 1.041424912903502, 1.042827972991933, 1.053302557606791, 1.015840821425795)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(4, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.031173921014486, 1.039124959442368, 1.040559866856918, 1.051079190099108, 1.034040096803066, 1.035631214662186, 1.041549442060937, 1.042980828166868, 1.053474652145984, 1.015885409612717)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(5, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.031212952390629, 1.039154607468478, 1.040594280125393, 1.051116742727418, 1.034045924047822, 1.035653794135171, 1.041570341295524, 1.043006494426994, 1.053503547581758, 1.015892894986566)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(6, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.030944602610908, 1.038950741017918, 1.040357709023066, 1.050858576789601, 1.03400570947641, 1.035498518206026, 1.041426577521094, 1.042830015368595, 1.053304857140912, 1.015841417294008)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(7, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.029823885385196, 1.038098600131795, 1.039370413643726, 1.049780820156362, 1.033834004790019, 1.034849120312378, 1.040824266256492, 1.042092546066169, 1.052474376070255, 1.015626081519176)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(8, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.02785269052195, 1.036597097290524, 1.037636574275911, 1.047886838740971, 1.033517876047131, 1.033703490374774, 1.039757728489602, 1.040793822564146, 1.05101112154261, 1.015246031722789)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(9, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.026541109938047, 1.035596260364477, 1.036484758559639, 1.046627774527953, 1.033298059274993, 1.032938935212846, 1.039043302073461, 1.039928626042457, 1.050035826165139, 1.014992288706293)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(10, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.025973805318399, 1.035162949645087, 1.035986999735316, 1.046083462219734, 1.033200737683691, 1.032607700421033, 1.038733159061709, 1.039554151485471, 1.049613582852464, 1.014882331617099)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(11, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.025763177590587, 1.035002009398828, 1.035802259786562, 1.045881413223575, 1.033164267567382, 1.032484639718184, 1.038617840650389, 1.039415080648025, 1.049456754494453, 1.014841476392688)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(12, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.025808353653391, 1.03503653111581, 1.035841880322257, 1.04592474737617, 1.033172105022205, 1.032511037783227, 1.038642582134412, 1.039444910628962, 1.04949039415722, 1.01485024052472)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(13, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.025956392828617, 1.035149646045988, 1.035971726006073, 1.04606675806542, 1.033197729586379, 1.032597528694671, 1.038723629194202, 1.039542655314234, 1.049600619124901, 1.014878954759077)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(14, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.026047617214717, 1.035219341355202, 1.036051748128263, 1.046154273199146, 1.033213475277388, 1.032650815278894, 1.038773549437657, 1.039602882488152, 1.04966853395409, 1.014896644934907)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(15, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.026578773239622, 1.035625019159379, 1.03651781408471, 1.046663917261075, 1.033304473188063, 1.032960914548785, 1.039063868688288, 1.039953482207588, 1.050063850598528, 1.014999584450991)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(16, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.02691212005, 1.03587950706847, 1.036810429967395, 1.046983838070717, 1.033360981602661, 1.033155385057744, 1.039245767359994, 1.040173448351838, 1.050311841155814, 1.015064133339927)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(17, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.027106615169807, 1.03602795105632, 1.036981202773168, 1.047170526352409, 1.033393735516866, 1.033268799216087, 1.039351789297158, 1.040301766342431, 1.050456496109328, 1.015101775412913)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(18, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.027172943035404, 1.036078567494564, 1.037039447928872, 1.047234196399547, 1.033404868704139, 1.033307467556478, 1.039387927021246, 1.040345522038992, 1.050505820720089, 1.015114608981917)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(19, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.026876348919662, 1.035852202342101, 1.036779025221838, 1.046949504904342, 1.033354940136152, 1.033134521982237, 1.039226259254243, 1.040149846456616, 1.050285233456556, 1.015057208700073)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(20, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.025912796361047, 1.03511633620926, 1.03593348553725, 1.046024935774131, 1.033190192632376, 1.032572059966474, 1.038699766086753, 1.039513871219856, 1.049568160285203, 1.014870499467663)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(21, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.025307519951773, 1.034653730503314, 1.035402729030267, 1.045444391970027, 1.033084755205689, 1.032218271726618, 1.03836805990252, 1.039114157681141, 1.04911737591888, 1.014753037223191)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(22, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.025628336028939, 1.03489895998277, 1.035684010215139, 1.045752075617657, 1.033140825035274, 1.03240583497896, 1.038543967508612, 1.039326038748887, 1.049356338232401, 1.014815312734219)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(23, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.026892512169831, 1.035864540155871, 1.036793215381517, 1.046965018325417, 1.03335767065292, 1.033143949161308, 1.03923507436399, 1.040160511089536, 1.050297256311398, 1.01506033767215)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(24, $colIndexes[$i]).Value = $rowValues[$i]
}

$rowValues = @(1.02, 1.028361849177327, 1.036985250885567, 1.038084102679324, 1.04837585421227, 1.033601204542138, 1.033999810167349, 1.040034059915514, 1.04112947093697, 1.051389378658239, 1.015344352060571)
for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item(25, $colIndexes[$i]).Value = $rowValues[$i]
}

Write-Host "Updated rows 2-25 with new vm_pu values (380 kV case, B column 1.05 -> 1.02)"